$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.760.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "'2.617.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'593.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'155.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.549"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").Value = "'2.615.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +9.91%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "'27.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").Value = "'3.077.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'67.647.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "'2.616.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "'11.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "'364.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "'4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  -4.27%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'67.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'9.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.52%  "
$ws.Range("D27").Value = "'2.743.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").Value = "'0.0000104"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "'577.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.09%  "
$ws.Range("D30").Value = "'1.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.74%  "
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("D32").Value = "'7.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'0.133"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  -4.15%  "
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("D38").Value = "'158.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "'0.370"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").Value = "'5.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'2.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.96%  "
$ws.Range("D44").Value = "'41.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'16.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'155.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").Value = "'0.0₆0287"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.59%  "
$ws.Range("D49").Value = "'3.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "'0.627"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "'20.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.49%  "
